# 10/21/2020: Updated QA automation standard guidelines & comments in
# FDR-Validation-Voucher and Voucher package.
#
# - Bump the two "Global Variables" dates (B1/B2) forward by 9 days; every
#   Voucher sheet pulls these via ='Global Variables'!B1 / !B2, so the
#   cached formula results on Voucher1-5 cascade automatically on recalc.
# - Move the active selection on "Global Variables" to B3.
# - Move the active tab from Voucher5 to Voucher1.

$wb = $excel.ActiveWorkbook

# --- Global Variables: refresh the two driving dates -----------------
$wsGlobal = $wb.Worksheets.Item("Global Variables")
$wsGlobal.Range("B1").Value = 44111
$wsGlobal.Range("B2").Value = 44113

# Record the new selection on this sheet (Excel only persists a sheet's
# selection while that sheet is active, so activate it first).
$wsGlobal.Activate()
$wsGlobal.Range("B3").Select() | Out-Null

# --- Switch the active/visible tab from Voucher5 to Voucher1 ---------
$wsVoucher1 = $wb.Worksheets.Item("Voucher1")
$wsVoucher1.Activate()
